$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: 57÷6= -> 89÷5=
$t.Cell(1, 1).Range.Text = "89÷5="
# Row 1, Col 2: 32÷3= -> 33÷4=
$t.Cell(1, 2).Range.Text = "33÷4="
# Row 1, Col 3: 66÷8= -> 96÷8=
$t.Cell(1, 3).Range.Text = "96÷8="
# Row 1, Col 4: 35÷6= -> 74÷3=
$t.Cell(1, 4).Range.Text = "74÷3="
# Row 1, Col 5: 77÷7= -> 50÷8=
$t.Cell(1, 5).Range.Text = "50÷8="

# Row 5, Col 1: 46÷9= -> 89÷4=
$t.Cell(5, 1).Range.Text = "89÷4="
# Row 5, Col 2: 42÷9= -> 44÷4=
$t.Cell(5, 2).Range.Text = "44÷4="
# Row 5, Col 3: 87÷6= -> 22÷3=
$t.Cell(5, 3).Range.Text = "22÷3="
# Row 5, Col 4: 41÷7= -> 31÷6=
$t.Cell(5, 4).Range.Text = "31÷6="
# Row 5, Col 5: 64÷7= -> 60÷9=
$t.Cell(5, 5).Range.Text = "60÷9="

# Row 9, Col 1: 28÷7= -> 10÷3=
$t.Cell(9, 1).Range.Text = "10÷3="
# Row 9, Col 2: 22÷9= -> 53÷8=
$t.Cell(9, 2).Range.Text = "53÷8="
# Row 9, Col 3: 12÷8= -> 24÷9=
$t.Cell(9, 3).Range.Text = "24÷9="
# Row 9, Col 4: 55÷8= -> 17÷3=
$t.Cell(9, 4).Range.Text = "17÷3="
# Row 9, Col 5: 37÷7= -> 92÷4=
$t.Cell(9, 5).Range.Text = "92÷4="

# Row 13, Col 1: 95÷9= -> 73÷6=
$t.Cell(13, 1).Range.Text = "73÷6="
# Row 13, Col 2: 83÷3= -> 94÷6=
$t.Cell(13, 2).Range.Text = "94÷6="
# Row 13, Col 3: 64÷7= -> 88÷8=
$t.Cell(13, 3).Range.Text = "88÷8="
# Row 13, Col 4: 29÷8= -> 32÷3=
$t.Cell(13, 4).Range.Text = "32÷3="
# Row 13, Col 5: 71÷8= -> 92÷8=
$t.Cell(13, 5).Range.Text = "92÷8="

# Row 17, Col 1: 35÷3= -> 93÷3=
$t.Cell(17, 1).Range.Text = "93÷3="
# Row 17, Col 2: 33÷7= -> 41÷9=
$t.Cell(17, 2).Range.Text = "41÷9="
# Row 17, Col 3: 58÷9= -> 30÷2=
$t.Cell(17, 3).Range.Text = "30÷2="
# Row 17, Col 4: 83÷8= -> 56÷6=
$t.Cell(17, 4).Range.Text = "56÷6="
# Row 17, Col 5: 53÷5= -> 99÷4=
$t.Cell(17, 5).Range.Text = "99÷4="

